$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: D1 "technology_abreviation" -> "abreviation" ---
$ws.Range("D1").Value = "abreviation"

# --- Forward-fill bus_name (col A) for every row within each bus block ---
$busGroups = @(
    @{Start=2; End=39; Value="b1"},
    @{Start=40; End=77; Value="b2"},
    @{Start=78; End=115; Value="b3"},
    @{Start=116; End=153; Value="b4"}
)
foreach ($g in $busGroups) {
    $ws.Range("A" + $g.Start + ":A" + $g.End).Value = $g.Value
}

# --- Forward-fill variable_name (col B) for every row within each variable block ---
$varGroups = @(
    @{Start=2; End=12; Value="PexistingR"},
    @{Start=13; End=19; Value="VexistingST"},
    @{Start=20; End=37; Value="PexistingCT"},
    @{Start=38; End=39; Value="PexistingG"},
    @{Start=40; End=50; Value="PexistingR"},
    @{Start=51; End=57; Value="VexistingST"},
    @{Start=58; End=75; Value="PexistingCT"},
    @{Start=76; End=77; Value="PexistingG"},
    @{Start=78; End=88; Value="PexistingR"},
    @{Start=89; End=95; Value="VexistingST"},
    @{Start=96; End=113; Value="PexistingCT"},
    @{Start=114; End=115; Value="PexistingG"},
    @{Start=116; End=126; Value="PexistingR"},
    @{Start=127; End=133; Value="VexistingST"},
    @{Start=134; End=151; Value="PexistingCT"},
    @{Start=152; End=153; Value="PexistingG"}
)
foreach ($g in $varGroups) {
    $ws.Range("B" + $g.Start + ":B" + $g.End).Value = $g.Value
}

# --- Column widths: best-fit for A, B, D; fixed width for C (technology_name) ---
$ws.Columns.Item(1).EntireColumn.AutoFit() | Out-Null
$ws.Columns.Item(2).EntireColumn.AutoFit() | Out-Null
$ws.Columns.Item(4).EntireColumn.AutoFit() | Out-Null
$ws.Columns.Item(3).ColumnWidth = 21.83

# --- Restore the authored selection/view ---
$ws.Range("N37").Select() | Out-Null
